$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet per the commit (timestamp moved from 0807_191502 to 0808_002655)
$ws.Name = "牌局記錄0808_002655"

# Append the 29 new game-history rows (840-868) that were recorded in this upload.
$newRows = @(
    @(6,"2,4",8,"K,3,5",2,"莊","紅","紅","藍","藍","藍","紅"),
    @(4,"8,4,2",6,"2,1,3",2,"莊","紅","紅","紅","藍","藍","藍"),
    @(0,"6,4,Q",4,"9,5",4,"莊","紅","紅","紅","藍","藍","藍"),
    @(3,"Q,3",9,"10,9",6,"莊","紅","紅","紅","藍","藍","藍"),
    @(0,"10,Q,K",2,"1,Q,1",2,"莊","紅","紅","紅","藍","藍","藍"),
    @(3,"J,3,J",2,"8,4,Q",1,"閒","紅","紅","紅","藍","藍","藍"),
    @(6,"8,8",8,"8,3,7",2,"莊","藍","紅","藍","紅","藍","紅"),
    @(8,"8,Q",1,"Q,1",7,"閒","藍","紅","藍","紅","藍","紅"),
    @(5,"K,J,5",3,"10,Q,3",2,"閒","藍","紅","藍","紅","藍","紅"),
    @(7,"4,3",7,"3,4",0,"和","紅","紅","藍","紅","紅","紅"),
    @(5,"5,10,J",7,"10,10,7",2,"莊","藍","藍","藍","紅","紅","紅"),
    @(6,"3,3",6,"10,Q,6",0,"和","紅","藍","藍","藍","紅","紅"),
    @(5,"1,4",9,"7,2",4,"莊","紅","藍","藍","藍","紅","紅"),
    @(5,"5,Q,K",6,"J,6",1,"莊","紅","藍","藍","藍","紅","紅"),
    @(0,"5,3,2",4,"K,K,4",4,"莊","藍","紅","紅","紅","藍","藍"),
    @(4,"4,Q,J",5,"9,6",1,"莊","紅","紅","紅","藍","藍","藍"),
    @(3,"6,9,8",6,"6,Q",3,"莊","紅","紅","紅","藍","藍","藍"),
    @(9,"1,4,4",7,"7,Q",2,"閒","紅","紅","紅","藍","藍","藍"),
    @(7,"10,7,J",0,"7,3,K",7,"閒","藍","藍","紅","紅","紅","藍"),
    @(4,"J,2,2",6,"9,3,4",2,"莊","藍","藍","紅","紅","紅","藍"),
    @(9,"4,5",2,"J,2",7,"閒","藍","藍","紅","紅","紅","藍"),
    @(9,"10,9",9,"5,4",0,"和","藍","紅","紅","紅","藍","藍"),
    @(5,"9,6,Q",3,"5,7,1",2,"閒","紅","紅","藍","藍","藍","紅"),
    @(9,"10,9",2,"Q,2",7,"閒","紅","紅","藍","藍","藍","紅"),
    @(7,"7,Q",1,"3,10,8",6,"閒","藍","藍","藍","紅","紅","紅"),
    @(6,"6,Q",5,"3,2,K",1,"閒","藍","藍","藍","紅","紅","紅"),
    @(7,"5,J,2",2,"8,4,Q",5,"閒","藍","藍","藍","紅","紅","紅"),
    @(0,"K,J",8,"8,K",8,"莊","藍","藍","紅","紅","紅","藍"),
    @(1,"1,10",9,"2,7",8,"莊","藍","藍","紅","紅","紅","藍")
)

$startRow = 840
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $newRows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le 12; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
